$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: id 103, email rock45@phptravels.com, password ramco, result na
$ws.Range("A4").Value = 103
$ws.Range("B4").Value = "rock45@phptravels.com"
$ws.Range("C4").Value = "ramco"
$ws.Range("D4").Value = "na"

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:rock45@phptravels.com")
$ws.Range("B4").Style = $ws.Range("B3").Style

# Row 5: id 104, email bunny@phptravels.com, password dessert, result na
$ws.Range("A5").Value = 104
$ws.Range("B5").Value = "bunny@phptravels.com"
$ws.Range("C5").Value = "dessert"
$ws.Range("D5").Value = "na"

$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:bunny@phptravels.com")
$ws.Range("B5").Style = $ws.Range("B3").Style

$null = $ws.Range("D4").Select()
